# Applies the 7.8 History Card & Advanced Story edit to 6_Mei.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update dialogue text in column B (History card rewrite) ---
$ws.Range("B7").Value = 'Every day during 5-7 PM, the Lord would rest and enjoy tea in the study next to the main hall until dinnertime.'
$ws.Range("B8").Value = 'Today, the Lord had a heated argument with Ming. After comforting Ming, I returned to the main hall and saw the Lord taking a nap.'
$ws.Range("B9").Value = 'Around 6 PM, I suddenly felt dizzy and nauseous, so I went to the study to inform the Lord that I wouldn’t be attending the banquet.'
$ws.Range("B10").Value = 'After that, I returned to my room and rested until Butler He came to inform me that something had happened to the Lord.'
$ws.Range("B11").Value = 'So your last encounter with the Lord was in the study?'
$ws.Range("B12").Value = 'Yes.'
$ws.Range("B13").Value = ' <color=#00CC00>(Cross-referencing with Butler He’s testimony, he last saw the Lord at the backyard entrance, which would have been after Mei''s visit.)</color>'
$ws.Range("B14").Value = ' <color=#00CC00>(His claim that Mei weren’t feeling well and that the Lord went to find the doctor is consistent.)</color>'
$ws.Range("B15").Value = 'Were you resting in your bed the entire time during the banquet?'
$ws.Range("B16").Value = 'Yes, I never left the room.'
$ws.Range("B17").Value = 'While resting, did you hear anything unusual?'
$ws.Range("B18").Value = 'I was bathing behind the screen at the time and didn’t hear anything unusual.'
$ws.Range("B19").Value = 'I only know that Ming came to visit me once.'
$ws.Range("B20").Value = 'But what’s strange is......I heard someone enter without knocking. I asked who it was, and Ming said it was him.'
$ws.Range("B21").Value = 'After a while, he seemed to close the door and leave.'
$ws.Range("B22").Value = 'You mean, you only heard Ming’s voice but didn’t actually see him come in?'
$ws.Range("B23").Value = 'Correct.'
$ws.Range("B24").Value = 'Do you remember what time Ming returned to your room?'
$ws.Range("B25").Value = 'I’m sorry......I truly can’t recall.'

# --- Row height corrections (text reflow changed wrapped-line count) ---
$ws.Rows.Item(7).RowHeight = 34
$ws.Rows.Item(14).RowHeight = 51

# --- Action2 (J) column: mark rows 10 and 16 as "appearAt" ---
$ws.Range("J10").Value = "appearAt"
$ws.Range("J16").Value = "appearAt"

# --- Update active selection (no more frozen topLeftCell scroll) ---
$ws.Range("H24").Select() | Out-Null

